# Re-enable the "Estate Tax" RAD test rows (previously marked DoNotRun) and
# stamp them with fresh execution timestamps, matching a re-run of the
# Katalon test-data generator against this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: Existing Liability with Notice/Invoice Number / Estate Tax
$ws.Range("C35").Value = "Y"
$ws.Range("B35").Value = "Thu Feb 06 13:23:46 EST 2025"

# Row 46: New Tax Return Amount Due / Estate Tax
$ws.Range("C46").Value = "Y"
$ws.Range("B46").Value = "Thu Feb 06 13:24:02 EST 2025"

# Update the view state saved with the sheet: selection now spans the whole
# Execute column.
$ws.Range("C2:C54").Select()
